# Updated some demo reports
# - rename the jxls placeholder tokens used by the report template
#   (case-only change: ITEM_NAME -> item_name, VOLUME -> volume)
# - leave the selection on the cell the author ended up on (C6)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The two placeholder cells referenced by both the table and the chart.
$ws.Range("A5").Value = '${row.item_name}'
$ws.Range("B5").Value = '${row.volume}'

# Move the active selection to C6 (matches the saved sheetView state).
[void]$ws.Range("C6").Select()
